# Commit: "add functions for debuging and calibration."
#
# 1) Sheet2: recalibrate the "motor 2" leg length (C23) from 150 -> 1.
#    Every downstream formula on Sheet2 depends on C23 transitively, so the
#    whole cascade of cached <v> results is refreshed by the normal
#    Excel auto-recalc that follows this script.
# 2) Sheet2: update the view state (selection moved to C24, scrolled down).
# 3) Add a new "Sheet3" calibration/debug worksheet with a cal_pos table
#    (B = raw motor angle, C = calibrated position via a linear fit
#    C = B*$C$1 + $D$1).

$wb = $excel.ActiveWorkbook

# --- 1) Sheet2: the actual data edit -------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("C23").Value = 1

# --- 2) Sheet2: view-state bookkeeping ------------------------------------
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("C24").Select()

# --- 3) New Sheet3: calibration / debug table -----------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)

$ws3.Range("A1").Value = "cal_pos"
$ws3.Range("C1").Value = 4.1611
$ws3.Range("D1").Value = 100

$motorAngles = @(20, 40, 60, 80, 90, 112, 130, 148, 166, 184, 202)
for ($i = 0; $i -lt $motorAngles.Length; $i++) {
    $row = $i + 2
    $ws3.Cells.Item($row, 2).Value = $motorAngles[$i]
}

# C2 typed individually, then C3:C12 filled from it (so C3:C12 share one
# formula group, matching how this was built by hand in Excel).
$ws3.Range("C2").Formula = "=SUM(B2*`$C`$1+`$D`$1)"
$ws3.Range("C3:C12").Formula = "=SUM(B3*`$C`$1+`$D`$1)"

$ws3.Range("C2:C12").NumberFormat = "0"
$ws3.Range("C6:C12").Select()
